$d = $word.ActiveDocument

# Replace a paragraph's visible text (excluding the paragraph mark) with new
# text. Used when the target text differs from the current text already, so
# a single assignment is enough to force the underlying run-merge / proofErr
# cleanup.
function Set-ParaTextDirect($idx, $text) {
    $para = $d.Paragraphs($idx)
    $rng = $para.Range
    $rng.MoveEnd(1, -1)
    $rng.Text = $text
}

# Same effect, but for paragraphs whose target text is identical (character
# for character) to what is already displayed. Assigning the same text as a
# no-op doesn't trigger the run-merge / proofErr cleanup, so first stash a
# different placeholder, then set the real text on the second pass.
function Set-ParaTextClean($idx, $text) {
    $para = $d.Paragraphs($idx)
    $rng = $para.Range
    $rng.MoveEnd(1, -1)
    $rng.Text = "TmpPlaceholder_" + $idx
    $para2 = $d.Paragraphs($idx)
    $rng2 = $para2.Range
    $rng2.MoveEnd(1, -1)
    $rng2.Text = $text
}

# For paragraphs where a stray proofErr marker sits at the very end (after
# the last run, touching the paragraph mark), plain Range.Text assignment
# leaves it behind. Deleting the whole paragraph (including its mark) and
# re-inserting clean text avoids the orphaned marker.
function Set-ParaTextViaDelete($idx, $text) {
    $para = $d.Paragraphs($idx)
    $rng = $d.Range($para.Range.Start, $para.Range.End)
    $rng.Delete()
    $rng.InsertBefore($text + [char]13)
}

# --- CREATE TABLE / INSERT block -------------------------------------------------
Set-ParaTextClean 3 "CREATE TABLE EMP1(emp_id int PRIMARY KEY,Basicsalary int,consalary int,Tax int);"
Set-ParaTextDirect 4 "INSERT INTO EMP1(emp_id,Basicsalary) VALUES(121,25000);"
Set-ParaTextDirect 5 "INSERT INTO EMP1(emp_id,Basicsalary) VALUES(122,30000);"
Set-ParaTextDirect 6 "INSERT INTO EMP1(emp_id,Basicsalary) VALUES(123,28000);"
Set-ParaTextDirect 7 "INSERT INTO EMP1(emp_id,Basicsalary) VALUES(124,10000);"
Set-ParaTextDirect 8 "INSERT INTO EMP1(emp_id,Basicsalary) VALUES(125,15000);"

# --- First PL/SQL block (TAX) -----------------------------------------------------
Set-ParaTextClean 11 "  CURSOR cemp IS SELECT * FROM EMP1;"
Set-ParaTextViaDelete 13 " FOR c1 in cemp"
Set-ParaTextClean 15 "  IF c1.BASICSALARY<20000 THEN"
Set-ParaTextClean 17 "  SET TAX=0.1*c1.BASICSALARY "
Set-ParaTextClean 19 "  ELSE IF c1.BASICSALARY>=20000 AND c1.BASICSALARY<30000 THEN"
Set-ParaTextClean 21 "  SET TAX=0.2*c1.BASICSALARY "
Set-ParaTextClean 25 "  SET TAX=0.3*c1.BASICSALARY "

# --- Second PL/SQL block (CONSALARY) ----------------------------------------------
Set-ParaTextClean 36 "  CURSOR cemp IS SELECT * FROM EMP1;"
Set-ParaTextViaDelete 38 " FOR c1 in cemp"
Set-ParaTextClean 40 "  IF c1.BASICSALARY<20000 THEN"
Set-ParaTextClean 42 "  SET CONSALARY=0.12*c1.BASICSALARY "
Set-ParaTextClean 44 "  ELSE IF c1.BASICSALARY>=20000 AND c1.BASICSALARY<30000 THEN"
Set-ParaTextClean 46 "  SET CONSALARY=0.16*c1.BASICSALARY "
Set-ParaTextClean 50 "  SET CONSALARY=0.21*c1.BASICSALARY "

# --- Third PL/SQL block (NETSALARY) -----------------------------------------------
Set-ParaTextClean 61 "  CURSOR cemp IS SELECT * FROM EMP1;"
Set-ParaTextViaDelete 63 " FOR c1 in cemp"
Set-ParaTextClean 65 "  IF c1.BASICSALARY<20000 THEN"
Set-ParaTextClean 69 "  ELSE IF c1.BASICSALARY>=20000 AND c1.BASICSALARY<30000 THEN"
